$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "25.758.04"
$ws.Range("E2").Value = "  +3.12%  "

Set-TextValue "D3" "1.676.93"
$ws.Range("E3").Value = "  +2.18%  "

Set-TextValue "D4" "0.9986"
$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue "D5" "237.22"
$ws.Range("E5").Value = "  +1.78%  "

Set-TextValue "D6" "0.9999"
$ws.Range("E6").Value = "  +0.00%  "

Set-TextValue "D7" "0.4627"
$ws.Range("E7").Value = "  -2.87%  "

Set-TextValue "D8" "0.2594"
$ws.Range("E8").Value = "  +0.24%  "

Set-TextValue "D9" "0.06140"
$ws.Range("E9").Value = "  +0.89%  "

Set-TextValue "D10" "1.672.13"
$ws.Range("E10").Value = "  +1.86%  "

Set-TextValue "D11" "0.06999"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("E12").Value = "  +2.40%  "

Set-TextValue "D13" "4.362"
$ws.Range("E13").Value = "  +0.81%  "

Set-TextValue "D14" "0.5764"
$ws.Range("E14").Value = "  -2.02%  "

Set-TextValue "D15" "75.37"
$ws.Range("E15").Value = "  +2.27%  "

Set-TextValue "D16" "0.9997"
$ws.Range("E16").Value = "  -0.01%  "

Set-TextValue "D17" "0.9999"
$ws.Range("E17").Value = "  +0.04%  "

Set-TextValue "D18" "25.754.25"
$ws.Range("E18").Value = "  +3.16%  "

Set-TextValue "D19" "0.000006697"
$ws.Range("E19").Value = "  +1.72%  "

Set-TextValue "D20" "11.42"
$ws.Range("E20").Value = "  +1.65%  "

Set-TextValue "D21" "1.885.28"
$ws.Range("E21").Value = "  +1.44%  "

Set-TextValue "D22" "4.464"
$ws.Range("E22").Value = "  +3.57%  "

$ws.Range("E23").Value = "  +1.37%  "

Set-TextValue "D24" "5.232"
$ws.Range("E24").Value = "  +0.12%  "

Set-TextValue "D25" "134.25"
$ws.Range("E25").Value = "  +0.55%  "

Set-TextValue "D26" "14.99"
$ws.Range("E26").Value = "  +0.58%  "

Set-TextValue "D27" "1.390"
$ws.Range("E27").Value = "  +0.78%  "

Set-TextValue "D28" "1.716"
$ws.Range("E28").Value = "  +4.62%  "

Set-TextValue "D29" "104.55"
$ws.Range("E29").Value = "  +0.34%  "

Set-TextValue "D30" "3.946"
$ws.Range("E30").Value = "  +1.68%  "

# Row 31/32 swap: row 31 (Filecoin) <-> row 32 (Stellar), with updated price/volume data
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D31" "0.07678"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D32" "3.613"
$ws.Range("E32").Value = "  +1.18%  "

Set-TextValue "D33" "0.04337"
$ws.Range("E33").Value = "  +1.54%  "

Set-TextValue "D34" "2.602"
$ws.Range("E34").Value = "  +1.25%  "

Set-TextValue "D35" "0.6087"
$ws.Range("E35").Value = "  +2.57%  "

Set-TextValue "D36" "0.9517"
$ws.Range("E36").Value = "  +2.37%  "

Set-TextValue "D37" "0.9338"
$ws.Range("E37").Value = "  +6.85%  "

Set-TextValue "D38" "109.22"
$ws.Range("E38").Value = "  +10.76%  "

Set-TextValue "D39" "2.446"
$ws.Range("E39").Value = "  -4.96%  "

Set-TextValue "D40" "0.9984"
$ws.Range("E40").Value = "  -0.10%  "

Set-TextValue "D41" "1.860"
$ws.Range("E41").Value = "  +5.17%  "

Set-TextValue "D42" "0.01448"
$ws.Range("E42").Value = "  -3.16%  "

Set-TextValue "D43" "5.059"
$ws.Range("E43").Value = "  +8.48%  "

Set-TextValue "D44" "0.3725"
$ws.Range("E44").Value = "  +0.43%  "

Set-TextValue "D45" "0.1117"
$ws.Range("E45").Value = "  +1.44%  "

Set-TextValue "D46" "0.05305"
$ws.Range("E46").Value = "  +1.99%  "

Set-TextValue "D47" "31.49"
$ws.Range("E47").Value = "  +9.83%  "

Set-TextValue "D48" "6.147"
$ws.Range("E48").Value = "  +0.60%  "

Set-TextValue "D49" "7.631"
$ws.Range("E49").Value = "  +6.95%  "

Set-TextValue "D50" "1.210"
$ws.Range("E50").Value = "  +2.28%  "

Set-TextValue "D51" "1.001"
$ws.Range("E51").Value = "  -0.01%  "
